{"js": "// The document contains a single table of simple arithmetic \"drill\" problems\n// (20 rows x 5 columns = 100 cells). The edit replaces the text of every\n// cell with a new problem, in row-major order, leaving the table's\n// structure (row/column count) untouched.\nconst newValues = [\"45-37=\", \"78+6=\", \"39+52=\", \"49+37=\", \"22-5=\", \"92-88=\", \"44-16=\", \"13+79=\", \"28+36=\", \"74-65=\", \"74-46=\", \"58+6=\", \"54-9=\", \"94-58=\", \"81-16=\", \"86+6=\", \"44-36=\", \"93-88=\", \"5+26=\", \"37+35=\", \"50-41=\", \"23+59=\", \"6+58=\", \"78-19=\", \"9+25=\", \"40-15=\", \"47+29=\", \"28+6=\", \"26+55=\", \"9+73=\", \"29+37=\", \"35-7=\", \"53-28=\", \"72-4=\", \"29+12=\", \"38+9=\", \"41-32=\", \"17+34=\", \"56-38=\", \"72-55=\", \"50-43=\", \"40-4=\", \"50-6=\", \"45-36=\", \"27+68=\", \"28+5=\", \"43-37=\", \"37+36=\", \"48+19=\", \"19+67=\", \"64-6=\", \"70-67=\", \"36+15=\", \"20-13=\", \"26+19=\", \"44+8=\", \"83-36=\", \"9+69=\", \"33+28=\", \"17+67=\", \"33+19=\", \"93-39=\", \"53+28=\", \"33-4=\", \"46+6=\", \"57+14=\", \"79+19=\", \"19+35=\", \"81-33=\", \"64-49=\", \"26-7=\", \"19+62=\", \"39+4=\", \"7+49=\", \"37+29=\", \"33+49=\", \"36+27=\", \"37+39=\", \"93-89=\", \"93-35=\", \"17+54=\", \"61-36=\", \"17+49=\", \"72-13=\", \"70-51=\", \"51-17=\", \"17+46=\", \"75-6=\", \"62-44=\", \"40-34=\", \"34+47=\", \"8+56=\", \"69+7=\", \"85-36=\", \"81-66=\", \"9+83=\", \"50-36=\", \"13+58=\", \"94-5=\", \"6+9=\"];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document but found none.\");\n}\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst rowCount = rows.items.length;\n// Load each row's cells so we know how many columns each row has.\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  const cells = rows.items[r].cells.items;\n  for (let c = 0; c < cells.length; c++) {\n    if (idx >= newValues.length) break;\n    const cell = table.getCell(r, c);\n    cell.value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single table of simple arithmetic \"drill\" problems\n# (20 rows x 5 columns = 100 cells). The edit replaces the text of every\n# cell with a new problem, in row-major order, leaving the table's\n# structure (row/column count) untouched.\n$newValues = @(\n  \"45-37=\", \"78+6=\", \"39+52=\", \"49+37=\", \"22-5=\",\n  \"92-88=\", \"44-16=\", \"13+79=\", \"28+36=\", \"74-65=\",\n  \"74-46=\", \"58+6=\", \"54-9=\", \"94-58=\", \"81-16=\",\n  \"86+6=\", \"44-36=\", \"93-88=\", \"5+26=\", \"37+35=\",\n  \"50-41=\", \"23+59=\", \"6+58=\", \"78-19=\", \"9+25=\",\n  \"40-15=\", \"47+29=\", \"28+6=\", \"26+55=\", \"9+73=\",\n  \"29+37=\", \"35-7=\", \"53-28=\", \"72-4=\", \"29+12=\",\n  \"38+9=\", \"41-32=\", \"17+34=\", \"56-38=\", \"72-55=\",\n  \"50-43=\", \"40-4=\", \"50-6=\", \"45-36=\", \"27+68=\",\n  \"28+5=\", \"43-37=\", \"37+36=\", \"48+19=\", \"19+67=\",\n  \"64-6=\", \"70-67=\", \"36+15=\", \"20-13=\", \"26+19=\",\n  \"44+8=\", \"83-36=\", \"9+69=\", \"33+28=\", \"17+67=\",\n  \"33+19=\", \"93-39=\", \"53+28=\", \"33-4=\", \"46+6=\",\n  \"57+14=\", \"79+19=\", \"19+35=\", \"81-33=\", \"64-49=\",\n  \"26-7=\", \"19+62=\", \"39+4=\", \"7+49=\", \"37+29=\",\n  \"33+49=\", \"36+27=\", \"37+39=\", \"93-89=\", \"93-35=\",\n  \"17+54=\", \"61-36=\", \"17+49=\", \"72-13=\", \"70-51=\",\n  \"51-17=\", \"17+46=\", \"75-6=\", \"62-44=\", \"40-34=\",\n  \"34+47=\", \"8+56=\", \"69+7=\", \"85-36=\", \"81-66=\",\n  \"9+83=\", \"50-36=\", \"13+58=\", \"94-5=\", \"6+9=\"\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($idx -ge $newValues.Length) { break }\n        $cell = $tbl.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
